$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..45 (Week, Day, Activity, Study Hours).
# Row 1 (header: Week/Day/Activity/Study Hours) is unchanged.
$data = @(
    @(3, 'Wednesday', 'Mat. Modelling - Midterm', 1),
    @(4, 'Tuesday', 'Mat. Modelling - Assignment 1 ', 3),
    @(4, 'Thursday', 'Mat. Modelling - Assignment 1 ', 1),
    @(4, 'Saturday', 'Logistics - 1-4th week review', 2),
    @(4, 'Saturday', 'S.C. Analytical Methods - Midterm', 3),
    @(4, 'Sunday', 'S.C. Analytical Methods - Midterm', 3),
    @(5, 'Saturday', 'Mat. Modelling - Midterm', 4),
    @(5, 'Saturday', 'S.C. Analytical Methods - Project Proposal ', 1),
    @(5, 'Sunday', 'S.C. Analytical Methods - Project Proposal ', 5),
    @(6, 'Wednesday', 'S.C. Analytical Methods - Project Proposal 2', 1),
    @(6, 'Saturday', 'Logistics - Midterm', 6),
    @(6, 'Sunday', 'Logistics - Midterm', 1),
    @(6, 'Sunday', 'Mat. Modelling - Project Report 1 ', 5),
    @(7, 'Monday', 'Logistics - 4-8th week review', 2),
    @(7, 'Monday', 'Mat. Modelling - Assignment 2', 1),
    @(7, 'Tuesday', 'Mat. Modelling - Assignment 2', 3),
    @(7, 'Wednesday', 'S.C. Analytical Methods - Final', 1),
    @(7, 'Thursday', 'Logistics - 4-8th week review', 2),
    @(7, 'Sunday', 'S.C. Analytical Methods - Project Proposal 2', 4),
    @(8, 'Wednesday', 'S.C. Analytical Methods - Project Proposal 2', 1),
    @(9, 'Monday', 'Mat. Modelling - Midterm', 1),
    @(9, 'Wednesday', 'Mat. Modelling - Project Report 2 & Presentation ', 1),
    @(9, 'Friday', 'S.C. Analytical Methods - Final', 1),
    @(9, 'Sunday', 'Logistics - 8-10th week review', 2),
    @(10, 'Wednesday', 'S.C. Analytical Methods - Project Proposal 2', 1),
    @(11, 'Wednesday', 'Mat. Modelling - Project Report 2 & Presentation ', 1),
    @(11, 'Saturday', 'Mat. Modelling - Final', 3),
    @(12, 'Monday', 'Mat. Modelling - Assignment 3', 3),
    @(12, 'Tuesday', 'Mat. Modelling - Project Report 2 & Presentation ', 3),
    @(12, 'Wednesday', 'S.C. Analytical Methods - Project Proposal 2', 1),
    @(12, 'Thursday', 'Logistics - 10-12th week review', 2),
    @(12, 'Thursday', 'Mat. Modelling - Assignment 3', 1),
    @(12, 'Friday', 'Mat. Modelling - Project Report 2 & Presentation ', 2),
    @(14, 'Monday', 'Mat. Modelling - Final', 3),
    @(14, 'Tuesday', 'Mat. Modelling - Final', 3),
    @(14, 'Wednesday', 'Mat. Modelling - Project Report 2 & Presentation ', 1),
    @(14, 'Thursday', 'Logistics - 12-14th week review', 2),
    @(14, 'Thursday', 'S.C. Analytical Methods - Final', 2),
    @(14, 'Friday', 'S.C. Analytical Methods - Final', 4),
    @(14, 'Sunday', 'S.C. Analytical Methods - Paper Review Presentation', 6),
    @(15, 'Monday', 'Logistics - Final', 3),
    @(15, 'Tuesday', 'Logistics - Final', 3),
    @(15, 'Wednesday', 'S.C. Analytical Methods - Final', 1),
    @(16, 'Monday', 'Logistics - Final', 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

Write-Host "Wrote $($data.Length) data rows (rows 2..$($data.Length + 1))"
